$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Hyperlinks.Delete()
$ws.Range("B2").ClearContents()
$ws.Range("C2").Clear()
$ws.Range("D2").ClearContents()
$ws.Range("C8").Select()
